# EE333-BOM.xlsx edit: begin KiCad Redesign sheet
# - Add a new "KiCad Redesign" worksheet after Sheet1, seeded with the same
#   header row (Component/Link/Cost Individual/Quantity/Cost/In kit?) as
#   Sheet1, bold like the Sheet1 header.
# - Make the new sheet the active tab, with the cursor left at C7.
# - Update Sheet1's selection to the header row (A1:F1) since it is no
#   longer the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Leave Sheet1's own selection on its header row before handing off the
# active tab to the new sheet.
$ws1.Activate() | Out-Null
$ws1.Range("A1:F1").Select() | Out-Null

# Add the new sheet right after Sheet1 and rename it.
$kicad = $wb.Worksheets.Add($null, $ws1)
$kicad.Name = "KiCad Redesign"

# Seed it with Sheet1's header row (same shared-string cells) and bold it
# to match the Sheet1 header formatting.
$kicad.Range("A1:F1").Value2 = $ws1.Range("A1:F1").Value2
$kicad.Range("A1:F1").Font.Bold = $true

# Approximate column widths for the new layout (engine quantizes to 1/6
# character steps, so these are the closest achievable values).
$ws1.Columns.Item(1).ColumnWidth = 14
$ws1.Columns.Item(2).ColumnWidth = 15.5
$ws1.Columns.Item(3).ColumnWidth = 12.666666666666666

# Make KiCad Redesign the active sheet/tab with the cursor at C7.
$kicad.Activate() | Out-Null
$kicad.Range("C7").Select() | Out-Null
